$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.324075666666667
$ws.Range("H2").Value = 3.972227
$ws.Range("I2").Value = 0.01675578032580584
$ws.Range("J2").Value = 0.01684165790066494
$ws.Range("M2").Value = 12.36292333333333
$ws.Range("N2").Value = 37.08877
$ws.Range("O2").Value = 0.918273862214392
$ws.Range("P2").Value = 0.9303794466068031
$ws.Range("Q2").Value = 16.36944595453222
$ws.Range("R2").Value = 147.32501359079
$ws.Range("S2").Value = 0.01538639511419365
$ws.Range("T2").Value = 0.01566913235756174
$ws.Range("G3").Value = 1.324075666666667
$ws.Range("H3").Value = 3.972227
$ws.Range("I3").Value = 0.01675578032580584
$ws.Range("J3").Value = 0.01684165790066494
$ws.Range("O3").Value = 0.04269179184247177
$ws.Range("P3").Value = 0.04325459680761149
$ws.Range("Q3").Value = 0.7610376468543333
$ws.Range("R3").Value = 6.849338821689001
$ws.Range("S3").Value = 0.0007153342858274865
$ws.Range("T3").Value = 0.0007284791220649865
$ws.Range("G4").Value = 1.324075666666667
$ws.Range("H4").Value = 3.972227
$ws.Range("I4").Value = 0.01675578032580584
$ws.Range("J4").Value = 0.01684165790066494
$ws.Range("M4").Value = 0.525528
$ws.Range("N4").Value = 1.051056
$ws.Range("O4").Value = 0.03903434594313629
$ws.Range("P4").Value = 0.02636595658558534
$ws.Range("Q4").Value = 0.695838836952
$ws.Range("R4").Value = 4.175033021712
$ws.Range("S4").Value = 0.000654050925784702
$ws.Range("T4").Value = 0.0004440464210382121
$ws.Range("I5").Value = 0.9679468703219594
$ws.Range("J5").Value = 0.9729078406975189
$ws.Range("M5").Value = 12.36292333333333
$ws.Range("N5").Value = 37.08877
$ws.Range("O5").Value = 0.918273862214392
$ws.Range("P5").Value = 0.9303794466068031
$ws.Range("Q5").Value = 945.6291305151077
$ws.Range("R5").Value = 8510.662174635969
$ws.Range("S5").Value = 0.8888403110288788
$ws.Range("T5").Value = 0.9051734584275773
$ws.Range("I6").Value = 0.9679468703219594
$ws.Range("J6").Value = 0.9729078406975189
$ws.Range("O6").Value = 0.04269179184247177
$ws.Range("P6").Value = 0.04325459680761149
$ws.Range("S6").Value = 0.0413233863023571
$ws.Range("T6").Value = 0.04208273638033509
$ws.Range("I7").Value = 0.9679468703219594
$ws.Range("J7").Value = 0.9729078406975189
$ws.Range("M7").Value = 0.525528
$ws.Range("N7").Value = 1.051056
$ws.Range("O7").Value = 0.03903434594313629
$ws.Range("P7").Value = 0.02636595658558534
$ws.Range("Q7").Value = 40.197174430536
$ws.Range("R7").Value = 241.183046583216
$ws.Range("S7").Value = 0.03778317299072344
$ws.Range("T7").Value = 0.02565164588960636
$ws.Range("G8").Value = 1.2088275
$ws.Range("H8").Value = 2.417655
$ws.Range("I8").Value = 0.0152973493522347
$ws.Range("J8").Value = 0.01025050140181618
$ws.Range("M8").Value = 12.36292333333333
$ws.Range("N8").Value = 37.08877
$ws.Range("O8").Value = 0.918273862214392
$ws.Range("P8").Value = 0.9303794466068031
$ws.Range("Q8").Value = 14.944641705725
$ws.Range("R8").Value = 89.66785023434998
$ws.Range("S8").Value = 0.01404715607131938
$ws.Range("T8").Value = 0.009536855821663999
$ws.Range("G9").Value = 1.2088275
$ws.Range("H9").Value = 2.417655
$ws.Range("I9").Value = 0.0152973493522347
$ws.Range("J9").Value = 0.01025050140181618
$ws.Range("O9").Value = 0.04269179184247177
$ws.Range("P9").Value = 0.04325459680761149
$ws.Range("Q9").Value = 0.6947965733474999
$ws.Range("R9").Value = 4.168779440085
$ws.Range("S9").Value = 0.000653071254287174
$ws.Range("T9").Value = 0.0004433813052114154
$ws.Range("G10").Value = 1.2088275
$ws.Range("H10").Value = 2.417655
$ws.Range("I10").Value = 0.0152973493522347
$ws.Range("J10").Value = 0.01025050140181618
$ws.Range("M10").Value = 0.525528
$ws.Range("N10").Value = 1.051056
$ws.Range("O10").Value = 0.03903434594313629
$ws.Range("P10").Value = 0.02636595658558534
$ws.Range("Q10").Value = 0.63527269842
$ws.Range("R10").Value = 2.54109079368
$ws.Range("S10").Value = 0.0005971220266281409
$ws.Range("T10").Value = 0.0002702642749407671
